# Applies scheduled market-data refresh values (Sheets/Halicarnassus_Profits.xlsx)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N17").Value = -7353.6921
$ws.Range("L17").Value = 7017.6921
$ws.Range("J17").Value = 2339.2307
$ws.Range("H17").Value = 2107.0967
$ws.Range("K19").Value = 108.77778
$ws.Range("M19").Value = 66.22221999999999
$ws.Range("I19").Value = 108.77778
$ws.Range("H19").Value = 290
$ws.Range("K40").Value = 5061.3335
$ws.Range("N40").Value = -6157.125
$ws.Range("M40").Value = -4886.3335
$ws.Range("L40").Value = 5807.125
$ws.Range("I40").Value = 5061.3335
$ws.Range("J40").Value = 5807.125
$ws.Range("H40").Value = 5412.294
$ws.Range("K74").Value = 6466.6665
$ws.Range("M74").Value = -5530.6665
$ws.Range("I74").Value = 6466.6665
$ws.Range("H74").Value = 6466.6665
$ws.Range("K77").Value = 32333.3325
$ws.Range("M77").Value = -27653.3325
$ws.Range("I77").Value = 6466.6665
$ws.Range("H77").Value = 6466.6665
$ws.Range("N88").Value = -1595
$ws.Range("L88").Value = 783
$ws.Range("J88").Value = 783
$ws.Range("H88").Value = 2177.6667
$ws.Range("N91").Value = -3591
$ws.Range("L91").Value = 783
$ws.Range("J91").Value = 783
$ws.Range("H91").Value = 2177.6667
$ws.Range("K100").Value = 3225.6667
$ws.Range("N100").Value = -5642
$ws.Range("M100").Value = -2684.6667
$ws.Range("L100").Value = 4560
$ws.Range("I100").Value = 3225.6667
$ws.Range("J100").Value = 4560
$ws.Range("H100").Value = 3759.4
$ws.Range("K106").Value = 2246.5
$ws.Range("M106").Value = -1615.5
$ws.Range("I106").Value = 2246.5
$ws.Range("H106").Value = 2331
$ws.Range("K111").Value = 2499.9999
$ws.Range("N111").Value = -9795.5
$ws.Range("M111").Value = 567.0001000000002
$ws.Range("L111").Value = 3661.5
$ws.Range("I111").Value = 833.3333
$ws.Range("J111").Value = 1220.5
$ws.Range("H111").Value = 988.2
$ws.Range("N112").Value = -10765.4999
$ws.Range("L112").Value = 8549.499899999999
$ws.Range("J112").Value = 2849.8333
$ws.Range("H112").Value = 2420
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
$ws.Range("I125").Value = 1000
$ws.Range("H125").Value = 1000
$ws.Range("K129").Value = 1652.25
$ws.Range("M129").Value = 3347.75
$ws.Range("I129").Value = 550.75
$ws.Range("H129").Value = 1217.909
$ws.Range("K131").Value = 6867.599999999999
$ws.Range("M131").Value = -1827.599999999999
$ws.Range("I131").Value = 2289.2
$ws.Range("H131").Value = 2037
$ws.Range("K132").Value = 33553.05
$ws.Range("M132").Value = -31023.05
$ws.Range("I132").Value = 11184.35
$ws.Range("H132").Value = 11184.35

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K12").Value = 999999
$ws.Range("M12").Value = -999826
$ws.Range("I12").Value = 999999
$ws.Range("H12").Value = 28500000
$ws.Range("K38").Value = 1987
$ws.Range("N38").Value = -5434
$ws.Range("M38").Value = -1520
$ws.Range("L38").Value = 4500
$ws.Range("I38").Value = 1987
$ws.Range("J38").Value = 4500
$ws.Range("H38").Value = 2489.6
$ws.Range("K97").Value = 810.38464
$ws.Range("M97").Value = -314.38464
$ws.Range("I97").Value = 810.38464
$ws.Range("H97").Value = 825.05884
$ws.Range("K102").Value = 2189.5789
$ws.Range("N102").Value = -9587.625
$ws.Range("M102").Value = -567.5789
$ws.Range("L102").Value = 6343.625
$ws.Range("I102").Value = 2189.5789
$ws.Range("J102").Value = 6343.625
$ws.Range("H102").Value = 3420.4075
$ws.Range("N105").Value = -31987.5
$ws.Range("L105").Value = 24999.5
$ws.Range("J105").Value = 24999.5
$ws.Range("H105").Value = 24999.5
$ws.Range("N128").Value = -129960
$ws.Range("L128").Value = 120000
$ws.Range("J128").Value = 120000
$ws.Range("H128").Value = 120000
$ws.Range("N133").Value = -185060
$ws.Range("L133").Value = 180000
$ws.Range("J133").Value = 180000
$ws.Range("H133").Value = 180000

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K16").Value = 654
$ws.Range("M16").Value = -484
$ws.Range("I16").Value = 654
$ws.Range("H16").Value = 654
$ws.Range("N64").Value = -1244.5
$ws.Range("L64").Value = 794.5
$ws.Range("J64").Value = 794.5
$ws.Range("H64").Value = 753.125
$ws.Range("N67").Value = -2354.5
$ws.Range("L67").Value = 794.5
$ws.Range("J67").Value = 794.5
$ws.Range("H67").Value = 753.125
$ws.Range("K80").Value = 745
$ws.Range("M80").Value = 253
$ws.Range("I80").Value = 745
$ws.Range("H80").Value = 1096.4445
$ws.Range("K83").Value = 3725
$ws.Range("M83").Value = 1267
$ws.Range("I83").Value = 745
$ws.Range("H83").Value = 1096.4445
$ws.Range("K86").Value = 3806
$ws.Range("M86").Value = -2683
$ws.Range("I86").Value = 3806
$ws.Range("H86").Value = 7830.4287
$ws.Range("K89").Value = 19030
$ws.Range("M89").Value = -13414
$ws.Range("I89").Value = 3806
$ws.Range("H89").Value = 7830.4287
$ws.Range("K99").Value = 781.6842
$ws.Range("M99").Value = 716.3158
$ws.Range("I99").Value = 781.6842
$ws.Range("H99").Value = 798.15
$ws.Range("K134").Value = 3088.3335
$ws.Range("M134").Value = -553.3335000000002
$ws.Range("I134").Value = 1029.4445
$ws.Range("H134").Value = 4031.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K99").Value = 2887.5
$ws.Range("M99").Value = -1389.5
$ws.Range("I99").Value = 2887.5
$ws.Range("H99").Value = 3189.2856
$ws.Range("K126").Value = 8662.5
$ws.Range("M126").Value = -6192.5
$ws.Range("I126").Value = 2887.5
$ws.Range("H126").Value = 3189.2856

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N39").Value = -18545.55
$ws.Range("L39").Value = 17957.55
$ws.Range("J39").Value = 5985.85
$ws.Range("H39").Value = 5985.85
$ws.Range("N55").Value = -7788
$ws.Range("L55").Value = 7434
$ws.Range("J55").Value = 2478
$ws.Range("H55").Value = 2015.6666
$ws.Range("K131").Value = 3414
$ws.Range("M131").Value = 1626
$ws.Range("I131").Value = 1138
$ws.Range("H131").Value = 1912.8
$ws.Range("N132").Value = -54571.25
$ws.Range("L132").Value = 49511.25
$ws.Range("J132").Value = 5501.25
$ws.Range("H132").Value = 5501.25
$ws.Range("N137").Value = -27471
$ws.Range("L137").Value = 17271
$ws.Range("J137").Value = 5757
$ws.Range("H137").Value = 4969.1665
$ws.Range("K139").Value = 1245
$ws.Range("M139").Value = 3895
$ws.Range("I139").Value = 415
$ws.Range("H139").Value = 415

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N64").ClearContents()
$ws.Range("L64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("L67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("K80").Value = 2569
$ws.Range("M80").Value = -1571
$ws.Range("I80").Value = 2569
$ws.Range("H80").Value = 2573.4285
$ws.Range("K83").Value = 12845
$ws.Range("M83").Value = -7853
$ws.Range("I83").Value = 2569
$ws.Range("H83").Value = 2573.4285
$ws.Range("K97").Value = 659.6667
$ws.Range("M97").Value = -163.6667
$ws.Range("I97").Value = 659.6667
$ws.Range("H97").Value = 763.58826
$ws.Range("N106").Value = -40524
$ws.Range("L106").Value = 38000
$ws.Range("J106").Value = 38000
$ws.Range("H106").Value = 38000
$ws.Range("K132").Value = 221523.87
$ws.Range("M132").Value = -218993.87
$ws.Range("I132").Value = 73841.28999999999
$ws.Range("H132").Value = 64854.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K82").Value = 1525
$ws.Range("M82").Value = -1164
$ws.Range("I82").Value = 1525
$ws.Range("H82").Value = 4915.875
$ws.Range("K85").Value = 1525
$ws.Range("M85").Value = -277
$ws.Range("I85").Value = 1525
$ws.Range("H85").Value = 4915.875
$ws.Range("N112").ClearContents()
$ws.Range("L112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("H112").Value = 0

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K81").Value = 2002
$ws.Range("M81").Value = -941
$ws.Range("I81").Value = 1001
$ws.Range("H81").Value = 1001
$ws.Range("K84").Value = 10010
$ws.Range("M84").Value = -4706
$ws.Range("I84").Value = 1001
$ws.Range("H84").Value = 1001
$ws.Range("N104").Value = -31430.334
$ws.Range("L104").Value = 24442.334
$ws.Range("J104").Value = 24442.334
$ws.Range("H104").Value = 24442.334
$ws.Range("N129").Value = -129999.5
$ws.Range("L129").Value = 119999.5
$ws.Range("J129").Value = 119999.5
$ws.Range("H129").Value = 119999.5
